$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# 1) "...tools such as Google Maps JavaScript API..." -> insert "the " before "Google Maps JavaScript API"
Replace-InParagraph 3 "tools such as Google Maps JavaScript API" "tools such as the Google Maps JavaScript API"

# 2) "Recent Earthquakes in the United States:" -> "30 Days of Earthquakes in the United States:"
Replace-InParagraph 5 "Recent Earthquakes in the United States:" "30 Days of Earthquakes in the United States:"

# 3) "San Francisco, Los Angeles, New York, among other cities" -> insert "and " before "New York"
Replace-InParagraph 20 "San Francisco, Los Angeles, New York, among other cities" "San Francisco, Los Angeles, and New York, among other cities"

# 3b) "...illustrates how the current version of each city..." -> "...illustrates how each city..."
Replace-InParagraph 20 "illustrates how the current version of each city" "illustrates how each city"

# 4) Conclusion paragraph: "demo! Please" -> "demo. Please"
Replace-InParagraph 29 "Thank you for viewing my demo! Please" "Thank you for viewing my demo. Please"

# 5) Final <h1> paragraph: "this demo! The CesiumJS" -> "this demo. The CesiumJS"
Replace-InParagraph 32 "create this demo! The CesiumJS" "create this demo. The CesiumJS"

# 6) Final <h1> paragraph: "(see TODO in repository)" -> "(see TODO items in repository)"
Replace-InParagraph 32 "see TODO in repository" "see TODO items in repository"
